$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.945.45"
$ws.Range("E2").Value = "  +7.45%  "
$ws.Range("D3").Value = "2.676.90"
$ws.Range("E3").Value = "  +11.33%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'511.98"
$ws.Range("E5").Value = "  +5.02%  "
$ws.Range("D6").Value = "'157.16"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.608"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").Value = "2.672.11"
$ws.Range("E9").Value = "  +10.36%  "
$ws.Range("E10").Value = "  +5.83%  "
$ws.Range("E11").Value = "  +5.14%  "
$ws.Range("E12").Value = "  +4.15%  "
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "3.138.65"
$ws.Range("E14").Value = "  +10.20%  "
$ws.Range("D15").Value = "60.967.46"
$ws.Range("E15").Value = "  +6.81%  "
$ws.Range("D16").Value = "'21.77"
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("E17").Value = "  +5.17%  "
$ws.Range("D18").Value = "2.673.83"
$ws.Range("E18").Value = "  +10.19%  "
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "'348.88"
$ws.Range("E20").Value = "  +7.53%  "
$ws.Range("D21").Value = "'10.51"
$ws.Range("E21").Value = "  +5.20%  "
$ws.Range("E22").Value = "  +3.54%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'60.30"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("D26").Value = "2.786.04"
$ws.Range("E26").Value = "  +9.81%  "
$ws.Range("E27").Value = "  +3.69%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "0.0₃0863"
$ws.Range("E29").Value = "  +10.46%  "
$ws.Range("D30").Value = "'7.55"
$ws.Range("E30").Value = "  +3.48%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "'157.51"
$ws.Range("E32").Value = "  +5.03%  "
$ws.Range("D33").Value = "'19.53"
$ws.Range("E33").Value = "  +5.29%  "
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D35").Value = "'5.72"
$ws.Range("E35").Value = "  +7.05%  "
$ws.Range("D36").Value = "'4.07"
$ws.Range("E36").Value = "  +9.19%  "
$ws.Range("E37").Value = "  +5.53%  "
$ws.Range("E38").Value = "  +11.22%  "
$ws.Range("D39").Value = "'310.37"
$ws.Range("E39").Value = "  +15.29%  "
$ws.Range("D40").Value = "'0.863"
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("E41").Value = "  +6.56%  "
$ws.Range("D42").Value = "'0.838"
$ws.Range("E42").Value = "  +29.04%  "
$ws.Range("D43").Value = "'35.43"
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("E44").Value = "  +8.74%  "
$ws.Range("D45").Value = "'0.0579"
$ws.Range("E45").Value = "  +8.82%  "
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "'0.999"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0237"
$ws.Range("E49").Value = "  +4.02%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'4.86"
$ws.Range("E50").Value = "  +6.15%  "
$ws.Range("D51").Value = "2.051.75"
$ws.Range("E51").Value = "  +9.97%  "
